$wb = $excel.ActiveWorkbook

# xlLineStyleNone
$xlNone = -4142
# RGB(0,0,0) - black
$black = 0
# Border edge indices: xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

function Set-TopBottomOnlyBorder($rng) {
    # Clear left/right first, then stamp top+bottom thin black.
    $rng.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
    $rng.Borders.Item($xlEdgeRight).LineStyle = $xlNone
    $rng.Borders.Item($xlEdgeTop).Color = $black
    $rng.Borders.Item($xlEdgeBottom).Color = $black
}

function Set-TopBottomRightBorder($rng) {
    # Clear left first, then stamp top+bottom+right thin black.
    $rng.Borders.Item($xlEdgeLeft).LineStyle = $xlNone
    $rng.Borders.Item($xlEdgeTop).Color = $black
    $rng.Borders.Item($xlEdgeBottom).Color = $black
    $rng.Borders.Item($xlEdgeRight).Color = $black
}

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomOnlyBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value2 = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomOnlyBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomOnlyBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value2 = "approach"
$ws2.Range("F2").Value2 = "approach"

# Remove the stray empty inline-string cell G5 entirely.
$ws2.Range("G5").ClearContents()
